$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily data for 15, 16, 17 marzo (rows 227-229)
$newRows = @(
    @{ Row = 227; A = 44301; B = 3; C = 32; D = 206.6782923206097 },
    @{ Row = 228; A = 44302; B = 7; C = 34; D = 219.5956855906478 },
    @{ Row = 229; A = 44303; B = 3; C = 32; D = 206.6782923206097 }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D

    # Copy the date-column style (border/bold/centered/date-format) from the
    # row above so the new column-A cells match the rest of the series.
    $ws.Cells.Item($r.Row - 1, 1).Copy()
    $ws.Cells.Item($r.Row, 1).PasteSpecial(-4122)
}

$excel.CutCopyMode = $false
